# Adding new RAAL model Production
# Updates GHI daily/hourly clear/cloudy sky figures with the new model's output.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Daily" sheet - row 2 (clear_sky_ghi/dni/dhi, cloudy_sky_ghi/dni/dhi)
# ---------------------------------------------------------------------------
$daily = $wb.Worksheets.Item("Daily")

$daily.Range("G2").Value = 2994.05
$daily.Range("H2").Value = 6264.4
$daily.Range("I2").Value = 734.29
$daily.Range("J2").Value = 775.17
$daily.Range("K2").Value = 19.89
$daily.Range("L2").Value = 772.89

# ---------------------------------------------------------------------------
# "Hourly" sheet - rows 9-19 (hour 7 through 17)
# Columns: H=clear_sky_ghi, I=clear_sky_dni, J=clear_sky_dhi,
#          K=cloudy_sky_ghi, L=cloudy_sky_dni, M=cloudy_sky_dhi
# ---------------------------------------------------------------------------
$hourly = $wb.Worksheets.Item("Hourly")

# Row 9 (hour 7)
$hourly.Range("I9").Value = 34.75
$hourly.Range("K9").Value = 3
$hourly.Range("M9").Value = 3

# Row 10 (hour 8)
$hourly.Range("I10").Value = 411.97
$hourly.Range("K10").Value = 49.32
$hourly.Range("L10").Value = 19.89
$hourly.Range("M10").Value = 47.05

# Row 11 (hour 9)
$hourly.Range("H11").Value = 243.02
$hourly.Range("I11").Value = 630.63
$hourly.Range("K11").Value = 60.75
$hourly.Range("L11").Value = 0
$hourly.Range("M11").Value = 60.75

# Row 12 (hour 10)
$hourly.Range("H12").Value = 368.65
$hourly.Range("I12").Value = 733.13
$hourly.Range("J12").Value = 85.47
$hourly.Range("K12").Value = 92.16
$hourly.Range("L12").Value = 0
$hourly.Range("M12").Value = 92.16

# Row 13 (hour 11)
$hourly.Range("H13").Value = 454.76
$hourly.Range("I13").Value = 784.24
$hourly.Range("K13").Value = 113.69
$hourly.Range("L13").Value = 0
$hourly.Range("M13").Value = 113.69

# Row 14 (hour 12)
$hourly.Range("H14").Value = 490.84
$hourly.Range("I14").Value = 802.66
$hourly.Range("J14").Value = 96.39
$hourly.Range("K14").Value = 122.71
$hourly.Range("L14").Value = 0
$hourly.Range("M14").Value = 122.71

# Row 15 (hour 13)
$hourly.Range("H15").Value = 472.89
$hourly.Range("I15").Value = 793.7
$hourly.Range("J15").Value = 94.90000000000001
$hourly.Range("K15").Value = 118.22
$hourly.Range("L15").Value = 0
$hourly.Range("M15").Value = 118.22

# Row 16 (hour 14)
$hourly.Range("H16").Value = 402.88
$hourly.Range("I16").Value = 754.83
$hourly.Range("J16").Value = 88.73999999999999
$hourly.Range("K16").Value = 100.72
$hourly.Range("L16").Value = 0
$hourly.Range("M16").Value = 100.72

# Row 17 (hour 15)
$hourly.Range("H17").Value = 289
$hourly.Range("I17").Value = 673.36
$hourly.Range("J17").Value = 76.95999999999999
$hourly.Range("K17").Value = 72.25
$hourly.Range("M17").Value = 72.25

# Row 18 (hour 16)
$hourly.Range("H18").Value = 147.45
$hourly.Range("I18").Value = 508.37
$hourly.Range("J18").Value = 56.73
$hourly.Range("K18").Value = 36.86
$hourly.Range("M18").Value = 36.86

# Row 19 (hour 17)
$hourly.Range("H19").Value = 21.86
$hourly.Range("I19").Value = 136.76

Write-Output "GHI production values updated"
